$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.104.05"
$ws.Range("E2").Value = "'  -1.10%  "
$ws.Range("D3").Value = "'3.172.68"
$ws.Range("E3").Value = "'  -4.62%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E5").Value = "'  -1.97%  "
$ws.Range("D6").Value = "'135.11"
$ws.Range("D8").Value = "'3.171.26"
$ws.Range("E8").Value = "'  -4.60%  "
$ws.Range("E9").Value = "'  -0.92%  "
$ws.Range("E10").Value = "'  -5.54%  "
$ws.Range("E11").Value = "'  -4.85%  "
$ws.Range("E12").Value = "'  -3.27%  "
$ws.Range("E13").Value = "'  -4.09%  "
$ws.Range("D14").Value = "'35.02"
$ws.Range("D15").Value = "'3.694.37"
$ws.Range("E15").Value = "'  -4.57%  "
$ws.Range("E16").Value = "'  -1.61%  "
$ws.Range("D17").Value = "'3.170.13"
$ws.Range("E17").Value = "'  -4.64%  "
$ws.Range("D18").Value = "'63.068.88"
$ws.Range("E18").Value = "'  -1.26%  "
$ws.Range("E19").Value = "'  -3.94%  "
$ws.Range("D20").Value = "'462.61"
$ws.Range("E20").Value = "'  -3.84%  "
$ws.Range("E21").Value = "'  -1.75%  "
$ws.Range("E22").Value = "'  -5.15%  "
$ws.Range("D23").Value = "'7.65"
$ws.Range("E23").Value = "'  -4.20%  "
$ws.Range("D24").Value = "'13.52"
$ws.Range("E24").Value = "'  -1.75%  "
$ws.Range("D25").Value = "'83.29"
$ws.Range("E25").Value = "'  -1.89%  "
$ws.Range("E26").Value = "'  -0.01%  "
$ws.Range("E27").Value = "'  +0.09%  "
$ws.Range("E28").Value = "'  -3.72%  "
$ws.Range("D29").Value = "'7.75"
$ws.Range("E29").Value = "'  -6.85%  "
$ws.Range("D30").Value = "'6.80"
$ws.Range("E30").Value = "'  -5.54%  "
$ws.Range("E31").Value = "'  -6.15%  "
$ws.Range("D32").Value = "'27.21"
$ws.Range("E32").Value = "'  -6.04%  "
$ws.Range("E33").Value = "'  -4.02%  "
$ws.Range("E34").Value = "'  -6.48%  "
$ws.Range("E35").Value = "'  -5.88%  "
$ws.Range("E36").Value = "'  -4.08%  "
$ws.Range("D37").Value = "'51.42"
$ws.Range("E37").Value = "'  -1.85%  "
$ws.Range("E38").Value = "'  -5.00%  "
$ws.Range("D39").Value = "'0.0390"
$ws.Range("E39").Value = "'  -3.03%  "
$ws.Range("D40").Value = "'404.61"
$ws.Range("E40").Value = "'  -6.94%  "
$ws.Range("D41").Value = "'8.14"
$ws.Range("E41").Value = "'  -2.53%  "
$ws.Range("E42").Value = "'  -5.12%  "
$ws.Range("D43").Value = "'2.815.24"
$ws.Range("E43").Value = "'  -9.13%  "
$ws.Range("D44").Value = "'2.61"
$ws.Range("E44").Value = "'  -6.33%  "
$ws.Range("E45").Value = "'  -5.90%  "
$ws.Range("E47").Value = "'  -6.56%  "
$ws.Range("D48").Value = "'25.41"
$ws.Range("E48").Value = "'  -4.16%  "
$ws.Range("D49").Value = "'123.97"
$ws.Range("E50").Value = "'  -1.68%  "
$ws.Range("D51").Value = "'34.23"
$ws.Range("E51").Value = "'  -7.55%  "
